# Lab 4 Soft engine
# Expands the final paragraph ("After created the second object ") into a
# run of new paragraphs describing modifiers / rendering, matching the
# target diff. The trailing "_GoBack" bookmark must end up attached to the
# very end of the (new) last paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 0: remove the _GoBack bookmark up front so none of the structural
# edits below have to fight with it; it gets re-created at the very end,
# once the final text is in place.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
}

# ---------------------------------------------------------------------
# Step 1: rewrite the text of the existing last paragraph.
# ---------------------------------------------------------------------
$rFirst = $d.Content
$null = $rFirst.Find.Execute(
    "After created the second object ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "After created the second object, you can use modifiers to make them attach to the previous item. This is known as solidifier.",
    2)

# ---------------------------------------------------------------------
# Step 2: append the skeleton of the new paragraphs as plain text, using
# "@@MARK@@" as a stand-in for a paragraph break (Find/Replace with a
# literal "^p" in the replacement text turns it into a real paragraph
# mark, and - importantly - carries the just-deleted bookmark's eventual
# re-insertion point to the end of the final paragraph). Paragraphs 4 and
# 5 (below) need extra inline markup (proofErr/lastRenderedPageBreak /
# multiple runs) so they are only seeded here as short placeholders and
# rebuilt via InsertXML afterwards.
# ---------------------------------------------------------------------
$skeletonParts = @(
    "Afterwards, if you want it to go outwards, you change the offset value from -1 to 1. This pushes it outwards.",
    "Afterwards, you can change the thickness level to the desired level using the thickness.",
    "@PARA4@",
    "@PARA5@",
    "",
    "Rendering the image:",
    "To do this, you click the render button, which brings you to the rendering engine. There are multiple slots which allows you to compare each engine. One could be better quality in terms of its calculations. However, that one would take longer to render."
)
$skeletonReplacement = "This is known as solidifier.@@MARK@@" + [string]::Join("@@MARK@@", $skeletonParts)

$rSkeleton = $d.Content
$null = $rSkeleton.Find.Execute(
    "This is known as solidifier.", $true, $false, $false, $false, $false,
    $true, 1, $false, $skeletonReplacement, 2)

$rMarks = $d.Content
$null = $rMarks.Find.Execute(
    "@@MARK@@", $true, $false, $false, $false, $false,
    $true, 1, $false, "^p", 2)

# ---------------------------------------------------------------------
# Step 3: replace the @PARA4@ placeholder paragraph with the fully
# formed paragraph (page-break hint + grammar-check proofErr markers
# around "an up and down buttons").
# ---------------------------------------------------------------------
$rPara4 = $d.Content
$null = $rPara4.Find.Execute("@PARA4@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rPara4.Text = ""
$xmlPara4 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">The order of the modifiers is from top to bottom. On each modifier, there are </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>an up and down buttons</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> which would allow you to change the order of the modifiers.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $rPara4.InsertXML($xmlPara4)

# ---------------------------------------------------------------------
# Step 4: replace the @PARA5@ placeholder paragraph with its two-run
# version.
# ---------------------------------------------------------------------
$rPara5 = $d.Content
$null = $rPara5.Find.Execute("@PARA5@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rPara5.Text = ""
$xmlPara5 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Another modifier is the array modifier which </w:t></w:r><w:r><w:t>allows you to create duplicates of the original object.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $rPara5.InsertXML($xmlPara5)

# ---------------------------------------------------------------------
# Step 5: re-create the _GoBack bookmark at the very end of the document
# (end of the last paragraph's run). Adding a bookmark exactly at the
# document's terminal offset is unreliable, so a scratch marker is
# appended first, the bookmark is anchored just before it, and the
# marker is removed again.
# ---------------------------------------------------------------------
$rEnd = $d.Content
$rEnd.Collapse(0)
$rEnd.InsertAfter("@@ENDMARK@@")

$rMarker = $d.Content
$null = $rMarker.Find.Execute("@@ENDMARK@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rBookmark = $rMarker.Duplicate
$rBookmark.Collapse(1)
$null = $d.Bookmarks.Add("_GoBack", $rBookmark)

$rMarker2 = $d.Content
$null = $rMarker2.Find.Execute("@@ENDMARK@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rMarker2.Delete()
